$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 91
$ws.Range("I2").Value = 270
$ws.Range("J2").Value = 1046
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 278
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = 178
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 103
$ws.Range("T2").Value = 197
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 1626
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1572
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 12
